# Add reference rows for the AXDW* (Data Warehouse) service accounts to the
# "Setup" sheet's account table (Table33), mirroring the existing AXDBAdmin
# row pattern (LEFT/MID around the environment code in $G$2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup")

# The account table ("Table33") currently covers A15:D22 - grow it by two
# rows so the new entries become part of the table (ref + autoFilter ref
# both expand).
$lo = $ws.ListObjects.Item("Table33")
$lo.Resize($ws.Range("A15:D24"))

# Row 23: AXDWAdmin
$ws.Range("A23").Value = "AXDWAdmin"
$ws.Range("B23").Formula = '=LEFT(A23,4) & $G$2 & MID(A23,5,10)'
$ws.Range("B23").NumberFormat = "General"
$ws.Range("C23").Formula = '= LEFT(A23,4) & $G$2 & MID(A23,5,10)'
$ws.Range("D23").Value = "SQL"

# Row 24: AXDWRuntimeuser
$ws.Range("A24").Value = "AXDWRuntimeuser"
$ws.Range("B24").Formula = '=LEFT(A24,4) & $G$2 & MID(A24,5,11)'
$ws.Range("B24").NumberFormat = "General"
$ws.Range("C24").Formula = '= LEFT(A24,4) & $G$2 & MID(A24,5,11)'
$ws.Range("D24").Value = "SQL"

# Move the cursor the way the authoring session left it.
$ws.Range("F32").Select() | Out-Null
